$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The spreadsheet computes "adjusted count" / "total cycles" from raw
# instruction / cache-access counts in column B; everything in columns
# D/E/F (and the summary cells F16/H16) is driven by formulas and will
# recalculate automatically once the inputs change.
$ws.Range("B6").Value = 1911   # ALU
$ws.Range("B7").Value = 162    # Jump
$ws.Range("B8").Value = 506    # Branch
$ws.Range("B9").Value = 365    # Memory
$ws.Range("B10").Value = 507   # Other
$ws.Range("B14").Value = 193   # Cache hit
$ws.Range("B15").Value = 264   # Cache miss

# Match the saved view/selection state.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F18").Select()
